$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 707.6667
$ws.Range("I6").Value = 438.57144
$ws.Range("K6").Value = 1315.71432
$ws.Range("M6").Value = -1203.71432

# ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 173.15
$ws.Range("I9").Value = 135.5
$ws.Range("J9").Value = 210.8
$ws.Range("K9").Value = 135.5
$ws.Range("L9").Value = 210.8
$ws.Range("M9").Value = 33.5
$ws.Range("N9").Value = -548.8

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1165693.9
$ws.Range("J17").Value = 1165693.9
$ws.Range("L17").Value = 3497081.7
$ws.Range("N17").Value = -3497417.7

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1122.683
$ws.Range("J112").Value = 1121.2821
$ws.Range("L112").Value = 3363.8463
$ws.Range("N112").Value = -5579.846299999999

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2879
$ws.Range("I138").Value = 1417.1538
$ws.Range("J138").Value = 4879.421
$ws.Range("K138").Value = 4251.4614
$ws.Range("L138").Value = 14638.263
$ws.Range("M138").Value = 888.5385999999999
$ws.Range("N138").Value = -24918.263

# ARM row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 29750
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 29750
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 29750
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -30268

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4858.026
$ws.Range("I32").Value = 4264
$ws.Range("K32").Value = 4264
$ws.Range("M32").Value = -3977

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1838.3334
$ws.Range("I45").Value = 1808.2
$ws.Range("J45").Value = 1989
$ws.Range("K45").Value = 1808.2
$ws.Range("L45").Value = 1989
$ws.Range("M45").Value = -1431.2
$ws.Range("N45").Value = -2743

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 775.0303
$ws.Range("I97").Value = 825.2143
$ws.Range("J97").Value = 494
$ws.Range("K97").Value = 825.2143
$ws.Range("L97").Value = 494
$ws.Range("M97").Value = -329.2143
$ws.Range("N97").Value = -1486

# ARM row 101
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 40999.668
$ws.Range("J101").Value = 40999.668
$ws.Range("L101").Value = 40999.668
$ws.Range("N101").Value = -47489.668

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1357.619
$ws.Range("I102").Value = 1184.7368
$ws.Range("K102").Value = 1184.7368
$ws.Range("M102").Value = 437.2632000000001

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2298.348
$ws.Range("I122").Value = 1842.3334
$ws.Range("K122").Value = 5527.0002
$ws.Range("M122").Value = -3077.0002

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 673.9091
$ws.Range("I94").Value = 632.9474
$ws.Range("J94").Value = 933.3333
$ws.Range("K94").Value = 632.9474
$ws.Range("L94").Value = 933.3333
$ws.Range("M94").Value = -181.9474
$ws.Range("N94").Value = -1835.3333

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1062.8
$ws.Range("I99").Value = 1003.2857
$ws.Range("J99").Value = 1201.6666
$ws.Range("K99").Value = 1003.2857
$ws.Range("L99").Value = 1201.6666
$ws.Range("M99").Value = 494.7143
$ws.Range("N99").Value = -4197.6666

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4880
$ws.Range("I105").Value = 5966.6665
$ws.Range("J105").Value = 3250
$ws.Range("K105").Value = 5966.6665
$ws.Range("L105").Value = 3250
$ws.Range("M105").Value = -4219.6665
$ws.Range("N105").Value = -6744

# BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 42750
$ws.Range("J126").Value = 42750
$ws.Range("L126").Value = 42750
$ws.Range("N126").Value = -52630

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1310.8524
$ws.Range("I134").Value = 969.56525
$ws.Range("J134").Value = 2357.4666
$ws.Range("K134").Value = 2908.69575
$ws.Range("L134").Value = 7072.399800000001
$ws.Range("M134").Value = -373.6957499999999
$ws.Range("N134").Value = -12142.3998

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1093.5
$ws.Range("I16").Value = 1167.75
$ws.Range("J16").Value = 945
$ws.Range("K16").Value = 1167.75
$ws.Range("L16").Value = 945
$ws.Range("M16").Value = -880.75
$ws.Range("N16").Value = -1519

# CRP row 93
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 19110.889
$ws.Range("I93").Value = 1519.6
$ws.Range("K93").Value = 1519.6
$ws.Range("M93").Value = 352.4000000000001

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1670
$ws.Range("I94").Value = 3294.5
$ws.Range("J94").Value = 1128.5
$ws.Range("K94").Value = 3294.5
$ws.Range("L94").Value = 1128.5
$ws.Range("M94").Value = -2843.5
$ws.Range("N94").Value = -2030.5

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1093.5
$ws.Range("I113").Value = 1167.75
$ws.Range("J113").Value = 945
$ws.Range("K113").Value = 1167.75
$ws.Range("L113").Value = 945
$ws.Range("M113").Value = 1002.25
$ws.Range("N113").Value = -5285

# CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2632.2
$ws.Range("I122").Value = 1770.6666
$ws.Range("J122").Value = 3924.5
$ws.Range("K122").Value = 5311.9998
$ws.Range("L122").Value = 11773.5
$ws.Range("M122").Value = -2861.9998
$ws.Range("N122").Value = -16673.5

# CUL row 19
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 2048.3333
$ws.Range("I19").Value = 290
$ws.Range("J19").Value = 2400
$ws.Range("K19").Value = 870
$ws.Range("L19").Value = 7200
$ws.Range("M19").Value = -696
$ws.Range("N19").Value = -7548

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 922.5333000000001
$ws.Range("J92").Value = 843.8182
$ws.Range("L92").Value = 2531.4546
$ws.Range("N92").Value = -5027.4546

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 6010.909
$ws.Range("I131").Value = 895
$ws.Range("J131").Value = 8934.286
$ws.Range("K131").Value = 2685
$ws.Range("L131").Value = 26802.858
$ws.Range("M131").Value = 2355
$ws.Range("N131").Value = -36882.858

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 973.75
$ws.Range("I102").Value = 971
$ws.Range("J102").Value = 987.5
$ws.Range("K102").Value = 971
$ws.Range("L102").Value = 987.5
$ws.Range("M102").Value = 651
$ws.Range("N102").Value = -4231.5

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1640.8823
$ws.Range("I113").Value = 1472.1333
$ws.Range("J113").Value = 2906.5
$ws.Range("K113").Value = 1472.1333
$ws.Range("L113").Value = 2906.5
$ws.Range("M113").Value = 697.8667
$ws.Range("N113").Value = -7246.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 145857.14
$ws.Range("I7").Value = 169333.33
$ws.Range("K7").Value = 169333.33
$ws.Range("M7").Value = -169221.33

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 77285.57000000001
$ws.Range("I40").Value = 104199.8
$ws.Range("K40").Value = 104199.8
$ws.Range("M40").Value = -104063.8

# LTW row 94
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 30330
$ws.Range("J94").Value = 30330
$ws.Range("L94").Value = 30330
$ws.Range("N94").Value = -31682

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 145857.14
$ws.Range("I126").Value = 169333.33
$ws.Range("K126").Value = 507999.99
$ws.Range("M126").Value = -505529.99

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8773821
$ws.Range("I136").Value = 1934.2903
$ws.Range("J136").Value = 47620748
$ws.Range("K136").Value = 5802.8709
$ws.Range("L136").Value = 142862244
$ws.Range("M136").Value = -3252.8709
$ws.Range("N136").Value = -142867344

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 7182.5625
$ws.Range("I100").Value = 276.25
$ws.Range("J100").Value = 27901.5
$ws.Range("K100").Value = 552.5
$ws.Range("L100").Value = 55803
$ws.Range("M100").Value = -11.5
$ws.Range("N100").Value = -56885

# WVR row 103
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 16500
$ws.Range("J103").Value = 16500
$ws.Range("L103").Value = 16500
$ws.Range("N103").Value = -18844

# WVR row 104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 2370
$ws.Range("J104").Value = 2370
$ws.Range("L104").Value = 2370
$ws.Range("N104").Value = -9358

# WVR row 105
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 51333.332
$ws.Range("J105").Value = 51333.332
$ws.Range("L105").Value = 51333.332
$ws.Range("N105").Value = -58321.332
